$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7466359478912636
$ws.Range("C2").Value = 0.2194227286667019
$ws.Range("E2").Value = 0.2261932789802152
$ws.Range("F2").Value = 2.334359484226638
$ws.Range("G2").Value = 0.002498382364234378
$ws.Range("J2").Value = 0.08082821524235229
$ws.Range("K2").Value = 0.3282684772701998
$ws.Range("L2").Value = 0.4073937232074485
$ws.Range("O2").Value = 4.546524183113121
# Row 3
$ws.Range("B3").Value = 0.7048634148921735
$ws.Range("C3").Value = 0.220401941611378
$ws.Range("E3").Value = 0.2249063562856648
$ws.Range("F3").Value = 2.336188907939999
$ws.Range("G3").Value = 0.002500910583379626
$ws.Range("J3").Value = 0.08016630187227491
$ws.Range("K3").Value = 0.2945379121028964
$ws.Range("L3").Value = 0.3977089723769893
$ws.Range("O3").Value = 4.579802438045277
# Row 4
$ws.Range("B4").Value = 0.6794496682884983
$ws.Range("C4").Value = 0.2210387314803395
$ws.Range("E4").Value = 0.2242013783275496
$ws.Range("F4").Value = 2.33843479985623
$ws.Range("G4").Value = 0.00250254568062742
$ws.Range("J4").Value = 0.0797619858871812
$ws.Range("K4").Value = 0.2738498744335232
$ws.Range("L4").Value = 0.3919205461470909
$ws.Range("O4").Value = 4.602429260024195
# Row 5
$ws.Range("B5").Value = 0.6691531091140064
$ws.Range("C5").Value = 0.2213071891178462
$ws.Range("E5").Value = 0.223935580976768
$ws.Range("F5").Value = 2.339632550437386
$ws.Range("G5").Value = 0.00250323286823763
$ws.Range("J5").Value = 0.07959776918221095
$ws.Range("K5").Value = 0.2654255216080088
$ws.Range("L5").Value = 0.3896016615478004
$ws.Range("O5").Value = 4.612201571455074
# Row 6
$ws.Range("B6").Value = 0.6674470031668989
$ws.Range("C6").Value = 0.2213523081362467
$ws.Range("E6").Value = 0.2238927451454522
$ws.Range("F6").Value = 2.339848507274453
$ws.Range("G6").Value = 0.002503348237678709
$ws.Range("J6").Value = 0.07957053455902141
$ws.Range("K6").Value = 0.2640270526403015
$ws.Range("L6").Value = 0.3892190311599677
$ws.Range("O6").Value = 4.613857578696894
# Row 7
$ws.Range("B7").Value = 0.6793105624166742
$ws.Range("C7").Value = 0.2210423156852315
$ws.Range("E7").Value = 0.2241977066077325
$ws.Range("F7").Value = 2.338449808864013
$ws.Range("G7").Value = 0.002502554863736623
$ws.Range("J7").Value = 0.0797597689730587
$ws.Range("K7").Value = 0.2737362349472505
$ws.Range("L7").Value = 0.3918891108478419
$ws.Range("O7").Value = 4.602558818985571
# Row 8
$ws.Range("B8").Value = 0.7321844698994084
$ws.Range("C8").Value = 0.2197529959778635
$ws.Range("E8").Value = 0.2257319022975928
$ws.Range("F8").Value = 2.334757408814838
$ws.Range("G8").Value = 0.002499236954439829
$ws.Range("J8").Value = 0.08059956124714418
$ws.Range("K8").Value = 0.3166337574390923
$ws.Range("L8").Value = 0.4040217358255376
$ws.Range("O8").Value = 4.557543216015702
# Row 9
$ws.Range("B9").Value = 0.8377069681289697
$ws.Range("C9").Value = 0.2175057378063272
$ws.Range("E9").Value = 0.2294140384173211
$ws.Range("F9").Value = 2.336415998212118
$ws.Range("G9").Value = 0.002493384397074278
$ws.Range("J9").Value = 0.08226240991939449
$ws.Range("K9").Value = 0.4009182028070768
$ws.Range("L9").Value = 0.4290604627650794
$ws.Range("O9").Value = 4.486674453546044
# Row 10
$ws.Range("B10").Value = 0.9163280136914409
$ws.Range("C10").Value = 0.2160246875901315
$ws.Range("E10").Value = 0.2325272036454678
$ws.Range("F10").Value = 2.343052591688121
$ws.Range("G10").Value = 0.002489479181443044
$ws.Range("J10").Value = 0.08349311606217213
$ws.Range("K10").Value = 0.4629240923294446
$ws.Range("L10").Value = 0.4482088192451101
$ws.Range("O10").Value = 4.445220121138334
# Row 11
$ws.Range("B11").Value = 0.9523272485369603
$ws.Range("C11").Value = 0.2153875506031753
$ws.Range("E11").Value = 0.2340314887944217
$ws.Range("F11").Value = 2.347246927105459
$ws.Range("G11").Value = 0.002487787453856625
$ws.Range("J11").Value = 0.08405480156245915
$ws.Range("K11").Value = 0.4911468225570275
$ws.Range("L11").Value = 0.4570818119353959
$ws.Range("O11").Value = 4.428666471843314
# Row 12
$ws.Range("B12").Value = 0.9659922868498825
$ws.Range("C12").Value = 0.2151515251239147
$ws.Range("E12").Value = 0.2346137336147933
$ws.Range("F12").Value = 2.349004046389553
$ws.Range("G12").Value = 0.002487158967766227
$ws.Range("J12").Value = 0.08426774521740299
$ws.Range("K12").Value = 0.5018359234073841
$ws.Range("L12").Value = 0.4604649505850205
$ws.Range("O12").Value = 4.422729404627546
# Row 13
$ws.Range("B13").Value = 0.9630478247327119
$ws.Range("C13").Value = 0.2152021244782603
$ws.Range("E13").Value = 0.2344877771755947
$ws.Range("F13").Value = 2.348618114289451
$ws.Range("G13").Value = 0.002487293784826862
$ws.Range("J13").Value = 0.08422187335238718
$ws.Range("K13").Value = 0.4995337641285857
$ws.Range("L13").Value = 0.4597353058383788
$ws.Range("O13").Value = 4.423993315638938
# Row 14
$ws.Range("B14").Value = 0.9534508242406901
$ws.Range("C14").Value = 0.2153680276416772
$ws.Range("E14").Value = 0.2340791381682692
$ws.Range("F14").Value = 2.347388104164366
$ws.Range("G14").Value = 0.002487735505188477
$ws.Range("J14").Value = 0.08407231574566509
$ws.Range("K14").Value = 0.4920261889438109
$ws.Range("L14").Value = 0.4573596824585024
$ws.Range("O14").Value = 4.428171382565523
# Row 15
$ws.Range("B15").Value = 0.9475766476603553
$ws.Range("C15").Value = 0.2154703305576291
$ws.Range("E15").Value = 0.2338304746201629
$ws.Range("F15").Value = 2.346656666448268
$ws.Range("G15").Value = 0.00248800765017264
$ws.Range("J15").Value = 0.08398073884423596
$ws.Range("K15").Value = 0.4874277966871716
$ws.Range("L15").Value = 0.4559075503506449
$ws.Range("O15").Value = 4.430773735766309
# Row 16
$ws.Range("B16").Value = 0.9139800254191357
$ws.Range("C16").Value = 0.2160670608846154
$ws.Range("E16").Value = 0.2324306623263297
$ws.Range("F16").Value = 2.342802115844023
$ws.Range("G16").Value = 0.002489591440791141
$ws.Range("J16").Value = 0.08345644391992479
$ws.Range("K16").Value = 0.4610799462046771
$ws.Range("L16").Value = 0.4476321963411607
$ws.Range("O16").Value = 4.446348315419868
# Row 17
$ws.Range("B17").Value = 0.8934290112823362
$ws.Range("C17").Value = 0.2164424964924017
$ws.Range("E17").Value = 0.2315944403672106
$ws.Range("F17").Value = 2.340738363090992
$ws.Range("G17").Value = 0.002490584717835075
$ws.Range("J17").Value = 0.08313526229052215
$ws.Range("K17").Value = 0.44492010037672
$ws.Range("L17").Value = 0.4425969557827472
$ws.Range("O17").Value = 4.45649303702794
# Row 18
$ws.Range("B18").Value = 0.8816306937512195
$ws.Range("C18").Value = 0.2166618832006542
$ws.Range("E18").Value = 0.2311217628680033
$ws.Range("F18").Value = 2.339661984994166
$ws.Range("G18").Value = 0.002491164007420346
$ws.Range("J18").Value = 0.08295070052041709
$ws.Range("K18").Value = 0.4356269161472142
$ws.Range("L18").Value = 0.4397161085067012
$ws.Range("O18").Value = 4.462544882074553
# Row 19
$ws.Range("B19").Value = 0.8776398019435305
$ws.Range("C19").Value = 0.2167367562084159
$ws.Range("E19").Value = 0.2309631490158921
$ws.Range("F19").Value = 2.33931654991548
$ws.Range("G19").Value = 0.002491361517681725
$ws.Range("J19").Value = 0.08288824140221251
$ws.Range("K19").Value = 0.4324806828051919
$ws.Range("L19").Value = 0.438743335609459
$ws.Range("O19").Value = 4.4646311770681
# Row 20
$ws.Range("B20").Value = 0.8956144218894053
$ws.Range("C20").Value = 0.2164021742099784
$ws.Range("E20").Value = 0.2316825995488472
$ws.Range("F20").Value = 2.34094660388395
$ws.Range("G20").Value = 0.002490478155961981
$ws.Range("J20").Value = 0.08316943480475558
$ws.Range("K20").Value = 0.4466401888531664
$ws.Range("L20").Value = 0.4431313848172209
$ws.Range("O20").Value = 4.455390668120032
# Row 21
$ws.Range("B21").Value = 0.9562688087949027
$ws.Range("C21").Value = 0.2153191556759708
$ws.Range("E21").Value = 0.2341988237620001
$ws.Range("F21").Value = 2.347744807896774
$ws.Range("G21").Value = 0.002487605432146836
$ws.Range("J21").Value = 0.08411623791758416
$ws.Range("K21").Value = 0.4942313020337963
$ws.Range("L21").Value = 0.458056834599887
$ws.Range("O21").Value = 4.426935186962766
# Row 22
$ws.Range("B22").Value = 0.9961014282332599
$ws.Range("C22").Value = 0.2146419007983695
$ws.Range("E22").Value = 0.2359167624321898
$ws.Range("F22").Value = 2.353171722633746
$ws.Range("G22").Value = 0.002485798644564981
$ws.Range("J22").Value = 0.08473645152537301
$ws.Range("K22").Value = 0.5253448932618596
$ws.Range("L22").Value = 0.4679461939791878
$ws.Range("O22").Value = 4.410269824669342
# Row 23
$ws.Range("B23").Value = 0.9748247396021839
$ws.Range("C23").Value = 0.2150005742768037
$ws.Range("E23").Value = 0.2349931670922842
$ws.Range("F23").Value = 2.35018531377969
$ws.Range("G23").Value = 0.002486756510196004
$ws.Range("J23").Value = 0.08440530749131625
$ws.Range("K23").Value = 0.5087382380868064
$ws.Range("L23").Value = 0.4626557999567069
$ws.Range("O23").Value = 4.41898763357392
# Row 24
$ws.Range("B24").Value = 0.8946263455726182
$ws.Range("C24").Value = 0.2164203928732853
$ws.Range("E24").Value = 0.2316427176176674
$ws.Range("F24").Value = 2.340852115249277
$ws.Range("G24").Value = 0.002490526306909232
$ws.Range("J24").Value = 0.08315398512946359
$ws.Range("K24").Value = 0.4458625451422336
$ws.Range("L24").Value = 0.4428897258916606
$ws.Range("O24").Value = 4.455888365342929
# Row 25
$ws.Range("B25").Value = 0.8089663984227116
$ws.Range("C25").Value = 0.2180837283585042
$ws.Range("E25").Value = 0.2283460942565938
$ws.Range("F25").Value = 2.335015315640831
$ws.Range("G25").Value = 0.002494898079276658
$ws.Range("J25").Value = 0.08181092866557194
$ws.Range("K25").Value = 0.3781013923133969
$ws.Range("L25").Value = 0.4221541453525788
$ws.Range("O25").Value = 4.503982687212186
